$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H ("Industries") values from row 34 through row 176 were switched
# from 1 to 0.
$ws.Range("H34:H176").Value = 0
